$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "Estado de Cuenta" detail rows (period / document / name / value reshuffled
# per the refreshed database export referenced in the commit message).

$ws.Range("C16").Value = "1047417915"
$ws.Range("D16").Value = "CARLOS ALFONSO PALENCIA RODRIGUEZ"
$ws.Range("E16").Value = "1806"
$ws.Range("G16").Value = 781242

$ws.Range("C17").Value = "1047417915"
$ws.Range("D17").Value = "CARLOS ALFONSO PALENCIA RODRIGUEZ"
$ws.Range("E17").Value = "1805"
$ws.Range("G17").Value = 1160000

$ws.Range("C18").Value = "1143363639"
$ws.Range("D18").Value = "OSCAR ALFONSO PALENCIA RODRIGUEZ"
$ws.Range("E18").Value = "1806"
$ws.Range("G18").Value = 1160000

$ws.Range("C19").Value = "1143363639"
$ws.Range("D19").Value = "OSCAR ALFONSO PALENCIA RODRIGUEZ"
$ws.Range("E19").Value = "1805"
$ws.Range("G19").Value = 781242

$ws.Range("C20").Value = "1047431310"
$ws.Range("D20").Value = "JORGE ELIECER PALENCIA RODRIGUEZ"
$ws.Range("E20").Value = "1806"
$ws.Range("G20").Value = 1160000

$ws.Range("C21").Value = "1047431310"
$ws.Range("D21").Value = "JORGE ELIECER PALENCIA RODRIGUEZ"
$ws.Range("E21").Value = "1805"
$ws.Range("G21").Value = 781242

$ws.Range("C22").Value = "1143353062"
$ws.Range("D22").Value = "HUGUER ENRIQUE PALENCIA RODRIGUEZ"
$ws.Range("E22").Value = "1806"
$ws.Range("G22").Value = 781242

$ws.Range("C23").Value = "1143353062"
$ws.Range("D23").Value = "HUGUER ENRIQUE PALENCIA RODRIGUEZ"
$ws.Range("E23").Value = "1805"
$ws.Range("G23").Value = 781242

$ws.Range("C24").Value = "1047416352"
$ws.Range("D24").Value = "JHON JAIRO PALENCIA RODRIGUEZ"
$ws.Range("E24").Value = "1806"
$ws.Range("G24").Value = 781242

$ws.Range("C25").Value = "1047416352"
$ws.Range("D25").Value = "JHON JAIRO PALENCIA RODRIGUEZ"
$ws.Range("E25").Value = "1805"
$ws.Range("G25").Value = 1160000

$ws.Range("C26").Value = "19935236"
$ws.Range("D26").Value = "EUGENIO VALDES MOSQUERA"
$ws.Range("E26").Value = "1806"
$ws.Range("G26").Value = 1160000

$ws.Range("C27").Value = "19935236"
$ws.Range("D27").Value = "EUGENIO VALDES MOSQUERA"
$ws.Range("E27").Value = "1805"
$ws.Range("G27").Value = 781242
